# Auto-generated cell updates applying the diff to cryptos.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '68.111.54'; ForceText = $false }
    @{ Cell = 'E2'; Value = '  -0.97%  '; ForceText = $false }
    @{ Cell = 'D3'; Value = '2.641.21'; ForceText = $false }
    @{ Cell = 'E4'; Value = '  -0.05%  '; ForceText = $false }
    @{ Cell = 'D5'; Value = '596.23'; ForceText = $true }
    @{ Cell = 'E5'; Value = '  -0.69%  '; ForceText = $false }
    @{ Cell = 'D6'; Value = '155.48'; ForceText = $true }
    @{ Cell = 'E6'; Value = '  -0.22%  '; ForceText = $false }
    @{ Cell = 'E7'; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = 'E8'; Value = '  -1.08%  '; ForceText = $false }
    @{ Cell = 'E9'; Value = '  -0.13%  '; ForceText = $false }
    @{ Cell = 'D10'; Value = '0.157'; ForceText = $true }
    @{ Cell = 'E10'; Value = '  -1.19%  '; ForceText = $false }
    @{ Cell = 'E11'; Value = '  -0.09%  '; ForceText = $false }
    @{ Cell = 'D12'; Value = '0.350'; ForceText = $true }
    @{ Cell = 'E12'; Value = '  -0.29%  '; ForceText = $false }
    @{ Cell = 'E13'; Value = '  -0.08%  '; ForceText = $false }
    @{ Cell = 'E14'; Value = '  -0.33%  '; ForceText = $false }
    @{ Cell = 'D15'; Value = '3.122.64'; ForceText = $false }
    @{ Cell = 'E15'; Value = '  -0.31%  '; ForceText = $false }
    @{ Cell = 'D16'; Value = '67.996.61'; ForceText = $false }
    @{ Cell = 'E16'; Value = '  -0.99%  '; ForceText = $false }
    @{ Cell = 'D17'; Value = '2.647.23'; ForceText = $false }
    @{ Cell = 'E17'; Value = '  -0.10%  '; ForceText = $false }
    @{ Cell = 'D18'; Value = '11.32'; ForceText = $true }
    @{ Cell = 'E18'; Value = '  -0.43%  '; ForceText = $false }
    @{ Cell = 'D19'; Value = '362.98'; ForceText = $true }
    @{ Cell = 'E19'; Value = '  -0.71%  '; ForceText = $false }
    @{ Cell = 'D20'; Value = '7.37'; ForceText = $true }
    @{ Cell = 'E20'; Value = '  -1.14%  '; ForceText = $false }
    @{ Cell = 'D21'; Value = '4.41'; ForceText = $true }
    @{ Cell = 'E21'; Value = '  +2.80%  '; ForceText = $false }
    @{ Cell = 'E22'; Value = '  -3.17%  '; ForceText = $false }
    @{ Cell = 'E23'; Value = '  -2.35%  '; ForceText = $false }
    @{ Cell = 'D24'; Value = '74.58'; ForceText = $true }
    @{ Cell = 'E24'; Value = '  +2.34%  '; ForceText = $false }
    @{ Cell = 'E25'; Value = '  +0.00%  '; ForceText = $false }
    @{ Cell = 'D26'; Value = '9.70'; ForceText = $true }
    @{ Cell = 'E26'; Value = '  -3.68%  '; ForceText = $false }
    @{ Cell = 'E27'; Value = '  -0.11%  '; ForceText = $false }
    @{ Cell = 'E28'; Value = '  -2.21%  '; ForceText = $false }
    @{ Cell = 'D29'; Value = '0.997'; ForceText = $true }
    @{ Cell = 'E29'; Value = '  -0.21%  '; ForceText = $false }
    @{ Cell = 'D30'; Value = '553.84'; ForceText = $true }
    @{ Cell = 'E30'; Value = '  -5.19%  '; ForceText = $false }
    @{ Cell = 'E31'; Value = '  -0.26%  '; ForceText = $false }
    @{ Cell = 'E32'; Value = '  -1.80%  '; ForceText = $false }
    @{ Cell = 'E33'; Value = '  -1.00%  '; ForceText = $false }
    @{ Cell = 'B34'; Value = 'FirstDigitalUSD'; ForceText = $false }
    @{ Cell = 'C34'; Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; ForceText = $false }
    @{ Cell = 'D34'; Value = '0.999'; ForceText = $true }
    @{ Cell = 'E34'; Value = '  -0.03%  '; ForceText = $false }
    @{ Cell = 'B35'; Value = 'Kaspa'; ForceText = $false }
    @{ Cell = 'C35'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; ForceText = $false }
    @{ Cell = 'D35'; Value = '0.128'; ForceText = $true }
    @{ Cell = 'E35'; Value = '  -2.48%  '; ForceText = $false }
    @{ Cell = 'E36'; Value = '  -0.82%  '; ForceText = $false }
    @{ Cell = 'D37'; Value = '161.13'; ForceText = $true }
    @{ Cell = 'E37'; Value = '  +0.54%  '; ForceText = $false }
    @{ Cell = 'E38'; Value = '  +0.29%  '; ForceText = $false }
    @{ Cell = 'E39'; Value = '  +0.91%  '; ForceText = $false }
    @{ Cell = 'E40'; Value = '  -3.48%  '; ForceText = $false }
    @{ Cell = 'D41'; Value = '5.30'; ForceText = $true }
    @{ Cell = 'E41'; Value = '  -1.61%  '; ForceText = $false }
    @{ Cell = 'D42'; Value = '0.0₆0334'; ForceText = $false }
    @{ Cell = 'E42'; Value = '  +3.76%  '; ForceText = $false }
    @{ Cell = 'E43'; Value = '  +0.35%  '; ForceText = $false }
    @{ Cell = 'E44'; Value = '  -2.26%  '; ForceText = $false }
    @{ Cell = 'E45'; Value = '  +0.00%  '; ForceText = $false }
    @{ Cell = 'D46'; Value = '159.38'; ForceText = $true }
    @{ Cell = 'E46'; Value = '  +2.01%  '; ForceText = $false }
    @{ Cell = 'E47'; Value = '  -0.60%  '; ForceText = $false }
    @{ Cell = 'D48'; Value = '21.98'; ForceText = $true }
    @{ Cell = 'E48'; Value = '  -0.49%  '; ForceText = $false }
    @{ Cell = 'E49'; Value = '  -1.59%  '; ForceText = $false }
    @{ Cell = 'E50'; Value = '  -0.14%  '; ForceText = $false }
    @{ Cell = 'E51'; Value = '  -0.79%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $target = $ws.Range($u.Cell)
    if ($u.ForceText) {
        $target.Value = "'" + $u.Value
    } else {
        $target.Value = $u.Value
    }
}

Write-Output "Applied $($updates.Count) cell updates"
